# Insert a new "Goals for this project" slide after the current slide 5
# ("How a Network Understands") and before the old slide 6
# ("FPGA Implementation"), using the "Two Content" layout, then fix a
# typo on the (now shifted) Device Utilization slide.

$p = $ppt.ActivePresentation

# "Two Content" is CustomLayout index 4 on this deck's single slide master.
$twoContentLayout = $p.SlideMaster.CustomLayouts.Item(4)

$newSlide = $p.Slides.AddSlide(6, $twoContentLayout)

# Title
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Goals for this project"

# Left content placeholder (sz="half" idx="1")
$left = $newSlide.Shapes.Item(2).TextFrame.TextRange
$left.Text = "Learn about FPGAs" + "`r`r`r`r" + "Implement accurate algorithms" + "`r`r" + "IEEE Floating Point"

# Right content placeholder (sz="half" idx="2")
$right = $newSlide.Shapes.Item(3).TextFrame.TextRange
$right.Text = "Implement a Fast Neural Network" + "`r`r" + "4 bit text recognition" + "`r`r`r" + "Scalable Network"

# Fix "Perf" -> "Performance" on the Device Utilization & Performance slide
# (was slide 13, now slide 14 after the insertion above).
$devSlide = $p.Slides.Item(14)
$titleRange = $devSlide.Shapes.Item(1).TextFrame.TextRange
$fullTitle = $titleRange.Text
$idx = $fullTitle.IndexOf("Perf")
if ($idx -ge 0) {
    $sub = $titleRange.Characters($idx + 1, 4)
    $sub.Text = "Performance"
}
